$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1667
$ws.Range("I40").Value = 1600.6
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 1600.6
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -1425.6
$ws.Range("N40").Value = -2100

$ws.Range("H116").Value = 10002039
$ws.Range("I116").Value = 22223776
$ws.Range("K116").Value = 22223776
$ws.Range("M116").Value = -22220334

$ws.Range("H137").Value = 3265.6494
$ws.Range("I137").Value = 1940.8846
$ws.Range("J137").Value = 3941.0195
$ws.Range("K137").Value = 5822.6538
$ws.Range("L137").Value = 11823.0585
$ws.Range("M137").Value = -3272.6538
$ws.Range("N137").Value = -16923.0585

$ws.Range("H138").Value = 3725.5737
$ws.Range("I138").Value = 1387.25
$ws.Range("J138").Value = 5242.324
$ws.Range("K138").Value = 4161.75
$ws.Range("L138").Value = 15726.972
$ws.Range("M138").Value = 978.25
$ws.Range("N138").Value = -26006.972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3678.913
$ws.Range("I61").Value = 4440.5
$ws.Range("J61").Value = 3272.7334
$ws.Range("K61").Value = 4440.5
$ws.Range("L61").Value = 3272.7334
$ws.Range("M61").Value = -4228.5
$ws.Range("N61").Value = -3696.7334

$ws.Range("H74").Value = 2555.5293
$ws.Range("I74").Value = 2858.0527
$ws.Range("J74").Value = 2172.3333
$ws.Range("K74").Value = 2858.0527
$ws.Range("L74").Value = 2172.3333
$ws.Range("M74").Value = -1984.0527
$ws.Range("N74").Value = -3920.3333

$ws.Range("H77").Value = 2555.5293
$ws.Range("I77").Value = 2858.0527
$ws.Range("J77").Value = 2172.3333
$ws.Range("K77").Value = 14290.2635
$ws.Range("L77").Value = 10861.6665
$ws.Range("M77").Value = -9922.263500000001
$ws.Range("N77").Value = -19597.6665

$ws.Range("H132").Value = 5419.5
$ws.Range("J132").Value = 3115.5652
$ws.Range("L132").Value = 9346.695599999999
$ws.Range("N132").Value = -14406.6956

$ws.Range("H136").Value = 3678.913
$ws.Range("I136").Value = 4440.5
$ws.Range("J136").Value = 3272.7334
$ws.Range("K136").Value = 13321.5
$ws.Range("L136").Value = 9818.200199999999
$ws.Range("M136").Value = -10771.5
$ws.Range("N136").Value = -14918.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 7264.4287
$ws.Range("I22").Value = 7264.4287
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 7264.4287
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -7091.4287
$ws.Range("N22").ClearContents()

$ws.Range("H134").Value = 2764.0278
$ws.Range("I134").Value = 2789.1667
$ws.Range("J134").Value = 2738.889
$ws.Range("K134").Value = 8367.500100000001
$ws.Range("L134").Value = 8216.667000000001
$ws.Range("M134").Value = -5832.500100000001
$ws.Range("N134").Value = -13286.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9999.6
$ws.Range("J4").Value = 9999.6
$ws.Range("L4").Value = 9999.6
$ws.Range("N4").Value = -10223.6

$ws.Range("H31").Value = 3618.9434
$ws.Range("I31").Value = 1457.1305
$ws.Range("J31").Value = 5276.3335
$ws.Range("K31").Value = 1457.1305
$ws.Range("L31").Value = 5276.3335
$ws.Range("M31").Value = -1162.1305
$ws.Range("N31").Value = -5866.3335

$ws.Range("H34").Value = 3618.9434
$ws.Range("I34").Value = 1457.1305
$ws.Range("J34").Value = 5276.3335
$ws.Range("K34").Value = 1457.1305
$ws.Range("L34").Value = 5276.3335
$ws.Range("M34").Value = -1255.1305
$ws.Range("N34").Value = -5680.3335

$ws.Range("H58").Value = 1801.9048
$ws.Range("I58").Value = 1538.8667
$ws.Range("J58").Value = 2459.5
$ws.Range("K58").Value = 1538.8667
$ws.Range("L58").Value = 2459.5
$ws.Range("M58").Value = -1335.8667
$ws.Range("N58").Value = -2865.5

$ws.Range("H122").Value = 1507.8422
$ws.Range("J122").Value = 1369
$ws.Range("L122").Value = 4107
$ws.Range("N122").Value = -9007

$ws.Range("H136").Value = 1801.9048
$ws.Range("I136").Value = 1538.8667
$ws.Range("J136").Value = 2459.5
$ws.Range("K136").Value = 4616.6001
$ws.Range("L136").Value = 7378.5
$ws.Range("M136").Value = -2066.6001
$ws.Range("N136").Value = -12478.5

$ws.Range("H137").Value = 78390
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 78390
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 78390
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -88590

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1805.1
$ws.Range("I5").Value = 2463.0833
$ws.Range("K5").Value = 7389.249899999999
$ws.Range("M5").Value = -7277.249899999999

$ws.Range("H122").Value = 952.4286
$ws.Range("I122").Value = 803.55554
$ws.Range("J122").Value = 1064.0834
$ws.Range("K122").Value = 7231.99986
$ws.Range("L122").Value = 9576.750599999999
$ws.Range("M122").Value = -4781.99986
$ws.Range("N122").Value = -14476.7506

$ws.Range("H132").Value = 1881
$ws.Range("I132").Value = 1444.4445
$ws.Range("J132").Value = 2238.182
$ws.Range("K132").Value = 13000.0005
$ws.Range("L132").Value = 20143.638
$ws.Range("M132").Value = -10470.0005
$ws.Range("N132").Value = -25203.638

$ws.Range("H135").Value = 1805.1
$ws.Range("I135").Value = 2463.0833
$ws.Range("K135").Value = 22167.7497
$ws.Range("M135").Value = -19632.7497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1720.8334
$ws.Range("I122").Value = 1065
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3195
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -745
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 3581
$ws.Range("I132").Value = 2845.6667
$ws.Range("K132").Value = 8537.000100000001
$ws.Range("M132").Value = -6007.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9401.143

$ws.Range("H22").Value = 1296.0952
$ws.Range("I22").Value = 2152.2
$ws.Range("J22").Value = 1028.5625
$ws.Range("K22").Value = 2152.2
$ws.Range("L22").Value = 1028.5625
$ws.Range("M22").Value = -1857.2
$ws.Range("N22").Value = -1618.5625

$ws.Range("H27").Value = 1296.0952
$ws.Range("I27").Value = 2152.2
$ws.Range("J27").Value = 1028.5625
$ws.Range("K27").Value = 2152.2
$ws.Range("L27").Value = 1028.5625
$ws.Range("M27").Value = -2045.2
$ws.Range("N27").Value = -1242.5625

$ws.Range("H82").Value = 1685.2413
$ws.Range("I82").Value = 1227.6666
$ws.Range("J82").Value = 2434
$ws.Range("K82").Value = 1227.6666
$ws.Range("L82").Value = 2434
$ws.Range("M82").Value = -866.6666
$ws.Range("N82").Value = -3156

$ws.Range("H85").Value = 1685.2413
$ws.Range("I85").Value = 1227.6666
$ws.Range("J85").Value = 2434
$ws.Range("K85").Value = 1227.6666
$ws.Range("L85").Value = 2434
$ws.Range("M85").Value = 20.33339999999998
$ws.Range("N85").Value = -4930

$ws.Range("H122").Value = 7696002
$ws.Range("I122").Value = 2049.5833
$ws.Range("J122").Value = 14290819
$ws.Range("K122").Value = 6148.749899999999
$ws.Range("L122").Value = 42872457
$ws.Range("M122").Value = -3698.749899999999
$ws.Range("N122").Value = -42877357

$ws.Range("H132").Value = 7575
$ws.Range("I132").Value = 8766.666999999999
$ws.Range("K132").Value = 26300.001
$ws.Range("M132").Value = -23770.001

$ws.Range("H136").Value = 2970.2
$ws.Range("I136").Value = 2243.1428
$ws.Range("K136").Value = 6729.428400000001
$ws.Range("M136").Value = -4179.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 11000000
$ws.Range("I4").Value = 11000000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 11000000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -10999887
$ws.Range("N4").ClearContents()

$ws.Range("H43").Value = 50000
$ws.Range("I43").Value = 50000
$ws.Range("K43").Value = 50000
$ws.Range("M43").Value = -49851

$ws.Range("H136").Value = 1981.7646
$ws.Range("I136").Value = 1670.25
$ws.Range("J136").Value = 2506.4211
$ws.Range("K136").Value = 5010.75
$ws.Range("L136").Value = 7519.263300000001
$ws.Range("M136").Value = -2460.75
$ws.Range("N136").Value = -12619.2633

Write-Output "applied 34 row edits across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
